# Add contributor info (name, email, repo link) as a second row, with
# email + repo link turned into hyperlinks, matching the "add my name and
# email and repo of the project in the excel sheet" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill C2/B2/A2 in this order so the shared-string table grows with the
# repo link first, then the email, then the name (matches target index order).
$ws.Range("C2").Value = "https://github.com/E0xMomen/open_source_project.git"
$ws.Range("B2").Value = "momenameer110@gmail.com"
$ws.Range("A2").Value = "Momen Ameer Abdelmomen Ali"

# Turn the email + repo link cells into real hyperlinks (adds the builtin
# "Hyperlink" style/font automatically).
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:momenameer110@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/E0xMomen/open_source_project.git")

# Widen the columns so the new content is fully visible.
$ws.Columns("A").ColumnWidth = 48.88
$ws.Columns("B").ColumnWidth = 41.02
$ws.Columns("C").ColumnWidth = 61.02

# Move the active selection to the newly filled-in cell.
$ws.Range("C2").Select() | Out-Null
